# Refresh the crypto price/volume snapshot (Price = column D, Volume(1h) = column E).
# Cells hold literal display text (not real numbers), e.g. "248.56" or "  +1.76%  ",
# mirroring the workbook as scraped. Where the new price text still parses as a plain
# number (no thousands separators), NumberFormat is forced to Text ("@") first so Excel
# keeps it as the exact literal string (with trailing zeros, etc.) instead of coercing it.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "30.542.58"
$ws.Range("E2").Value = "  +0.60%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "1.876.34"
$ws.Range("E3").Value = "  -0.10%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.09%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.56"
$ws.Range("E5").Value = "  +1.76%  "

# Row 6: USDC
$ws.Range("E6").Value = "  -0.05%  "

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4761"
$ws.Range("E7").Value = "  -0.08%  "

# Row 8: Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2919"
$ws.Range("E8").Value = "  +1.39%  "

# Row 9: Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06504"
$ws.Range("E9").Value = "  -0.35%  "

# Row 10: Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.01"
$ws.Range("E10").Value = "  +3.10%  "

# Row 11: TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07741"
$ws.Range("E11").Value = "  -0.31%  "

# Row 12: Polygon
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7387"
$ws.Range("E12").Value = "  -0.08%  "

# Row 13: Litecoin
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.62"
$ws.Range("E13").Value = "  +0.23%  "

# Row 14: WrappedEther
$ws.Range("D14").Value = "1.874.32"
$ws.Range("E14").Value = "  -0.30%  "

# Row 15: Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.181"
$ws.Range("E15").Value = "  +0.87%  "

# Row 16: BitcoinCash
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "273.99"
$ws.Range("E16").Value = "  -0.82%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "30.639.44"
$ws.Range("E17").Value = "  +0.93%  "

# Row 18: Avalanche
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.21"
$ws.Range("E18").Value = "  -1.31%  "

# Row 19: Dai
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.06%  "

# Row 20: ShibaInu
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007516"
$ws.Range("E20").Value = "  -0.36%  "

# Row 21: WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "2.119.83"
$ws.Range("E21").Value = "  -0.36%  "

# Row 22: BinanceUSD
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23: Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.244"
$ws.Range("E23").Value = "  +0.29%  "

# Row 24: Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.194"
$ws.Range("E24").Value = "  +0.38%  "

# Row 25: Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "165.49"
$ws.Range("E25").Value = "  +0.66%  "

# Row 26: Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.183"
$ws.Range("E26").Value = "  -0.56%  "

# Row 27: EthereumClassic
$ws.Range("E27").Value = "  -0.64%  "

# Row 28: LidoDAOToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.908"
$ws.Range("E28").Value = "  -2.38%  "

# Row 29: Stellar
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.09846"
$ws.Range("E29").Value = "  -1.13%  "

# Row 30: Toncoin
$ws.Range("E30").Value = "  -3.08%  "

# Row 31: PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.503"
$ws.Range("E31").Value = "  -0.64%  "

# Row 32: Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.267"
$ws.Range("E32").Value = "  -1.04%  "

# Row 33: InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.098"
$ws.Range("E33").Value = "  +0.29%  "

# Row 34: Hedera
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04814"
$ws.Range("E34").Value = "  +1.30%  "

# Row 35: ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.122"
$ws.Range("E35").Value = "  -0.03%  "

# Row 36: ImmutableX
$ws.Range("E36").Value = "  +0.05%  "

# Row 37: HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  -0.02%  "

# Row 38: VeChain
$ws.Range("E38").Value = "  +0.83%  "

# Row 39: MXToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.761"
$ws.Range("E39").Value = "  +0.40%  "

# Row 40: FraxShare
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.272"
$ws.Range("E40").Value = "  -0.42%  "

# Row 41: Aave
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.43"
$ws.Range("E41").Value = "  +5.91%  "

# Row 42: RenderToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.981"
$ws.Range("E42").Value = "  +3.78%  "

# Row 43: TheSandbox
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4207"
$ws.Range("E43").Value = "  +0.78%  "

# Row 44: PaxDollar
$ws.Range("E44").Value = "  +0.00%  "

# Row 45: TrustWalletToken
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8357"
$ws.Range("E45").Value = "  -0.82%  "

# Row 46: Quant
$ws.Range("E46").Value = "  +0.03%  "

# Row 47: EnergySwap
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.416"
$ws.Range("E47").Value = "  +2.08%  "

# Row 48: Aptos
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.984"
$ws.Range("E48").Value = "  -1.59%  "

# Row 49: Elrond
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.31"
$ws.Range("E49").Value = "  +0.36%  "

# Row 50: Maker
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "912.85"
$ws.Range("E50").Value = "  +0.21%  "

# Row 51: Cronos
$ws.Range("E51").Value = "  +1.59%  "
